$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 529.5
$ws.Range("I4").Value = 499.375
$ws.Range("J4").Value = 650
$ws.Range("K4").Value = 499.375
$ws.Range("L4").Value = 650
$ws.Range("M4").Value = -385.375
$ws.Range("N4").Value = -878
$ws.Range("H17").Value = 2369.5833
$ws.Range("J17").Value = 2369.5833
$ws.Range("L17").Value = 7108.749899999999
$ws.Range("N17").Value = -7444.749899999999
$ws.Range("H98").Value = 305.8387
$ws.Range("I98").Value = 307.17856
$ws.Range("J98").Value = 293.33334
$ws.Range("K98").Value = 307.17856
$ws.Range("L98").Value = 293.33334
$ws.Range("M98").Value = 1190.82144
$ws.Range("N98").Value = -3289.33334
$ws.Range("H122").Value = 305.8387
$ws.Range("I122").Value = 307.17856
$ws.Range("J122").Value = 293.33334
$ws.Range("K122").Value = 921.53568
$ws.Range("L122").Value = 880.0000200000001
$ws.Range("M122").Value = 1528.46432
$ws.Range("N122").Value = -5780.00002
$ws.Range("H132").Value = 199388
$ws.Range("I132").Value = 240690.53
$ws.Range("K132").Value = 722071.59
$ws.Range("M132").Value = -719541.59
$ws.Range("H137").Value = 529479.6
$ws.Range("I137").Value = 772321.9
$ws.Range("J137").Value = 3321.5
$ws.Range("K137").Value = 2316965.7
$ws.Range("L137").Value = 9964.5
$ws.Range("M137").Value = -2314415.7
$ws.Range("N137").Value = -15064.5
$ws.Range("H138").Value = 6367.224
$ws.Range("I138").Value = 3069.5
$ws.Range("J138").Value = 6576.603
$ws.Range("K138").Value = 9208.5
$ws.Range("L138").Value = 19729.809
$ws.Range("M138").Value = -4068.5
$ws.Range("N138").Value = -30009.809
$ws.Range("H141").Value = 3272.6965
$ws.Range("I141").Value = 2939.647
$ws.Range("K141").Value = 8818.940999999999
$ws.Range("M141").Value = -3638.940999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 30179.176
$ws.Range("I2").Value = 42361.25
$ws.Range("J2").Value = 942.2
$ws.Range("K2").Value = 42361.25
$ws.Range("L2").Value = 942.2
$ws.Range("M2").Value = -42248.25
$ws.Range("N2").Value = -1168.2
$ws.Range("H5").Value = 253.75
$ws.Range("I5").Value = 238.33333
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 238.33333
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -126.33333
$ws.Range("N5").Value = -524
$ws.Range("H32").Value = 23101.682
$ws.Range("I32").Value = 27132.719
$ws.Range("J32").Value = 14502.134
$ws.Range("K32").Value = 27132.719
$ws.Range("L32").Value = 14502.134
$ws.Range("M32").Value = -26845.719
$ws.Range("N32").Value = -15076.134
$ws.Range("H45").Value = 1823.963
$ws.Range("I45").Value = 1202.4
$ws.Range("J45").Value = 3599.8572
$ws.Range("K45").Value = 1202.4
$ws.Range("L45").Value = 3599.8572
$ws.Range("M45").Value = -825.4000000000001
$ws.Range("N45").Value = -4353.8572
$ws.Range("H61").Value = 5123.2
$ws.Range("I61").Value = 3464.9565
$ws.Range("J61").Value = 8301.5
$ws.Range("K61").Value = 3464.9565
$ws.Range("L61").Value = 8301.5
$ws.Range("M61").Value = -3252.9565
$ws.Range("N61").Value = -8725.5
$ws.Range("H74").Value = 7144276.5
$ws.Range("I74").Value = 8929417
$ws.Range("J74").Value = 3713.8572
$ws.Range("K74").Value = 8929417
$ws.Range("L74").Value = 3713.8572
$ws.Range("M74").Value = -8928543
$ws.Range("N74").Value = -5461.8572
$ws.Range("H77").Value = 7144276.5
$ws.Range("I77").Value = 8929417
$ws.Range("J77").Value = 3713.8572
$ws.Range("K77").Value = 44647085
$ws.Range("L77").Value = 18569.286
$ws.Range("M77").Value = -44642717
$ws.Range("N77").Value = -27305.286
$ws.Range("H116").Value = 30179.176
$ws.Range("I116").Value = 42361.25
$ws.Range("J116").Value = 942.2
$ws.Range("K116").Value = 42361.25
$ws.Range("L116").Value = 942.2
$ws.Range("M116").Value = -40067.25
$ws.Range("N116").Value = -5530.2
$ws.Range("H132").Value = 19396.473
$ws.Range("I132").Value = 27958.318
$ws.Range("J132").Value = 5942.143
$ws.Range("K132").Value = 83874.954
$ws.Range("L132").Value = 17826.429
$ws.Range("M132").Value = -81344.954
$ws.Range("N132").Value = -22886.429
$ws.Range("H136").Value = 5123.2
$ws.Range("I136").Value = 3464.9565
$ws.Range("J136").Value = 8301.5
$ws.Range("K136").Value = 10394.8695
$ws.Range("L136").Value = 24904.5
$ws.Range("M136").Value = -7844.869499999999
$ws.Range("N136").Value = -30004.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 30179.176
$ws.Range("I3").Value = 42361.25
$ws.Range("J3").Value = 942.2
$ws.Range("K3").Value = 42361.25
$ws.Range("L3").Value = 942.2
$ws.Range("M3").Value = -42247.25
$ws.Range("N3").Value = -1170.2
$ws.Range("H4").Value = 253.75
$ws.Range("I4").Value = 238.33333
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 238.33333
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -123.33333
$ws.Range("N4").Value = -530
$ws.Range("H20").Value = 2693.2
$ws.Range("I20").Value = 2166.75
$ws.Range("J20").Value = 4799
$ws.Range("K20").Value = 2166.75
$ws.Range("L20").Value = 4799
$ws.Range("M20").Value = -1919.75
$ws.Range("N20").Value = -5293
$ws.Range("H107").Value = 4920.364
$ws.Range("I107").Value = 3016
$ws.Range("J107").Value = 9998.666999999999
$ws.Range("K107").Value = 3016
$ws.Range("L107").Value = 9998.666999999999
$ws.Range("M107").Value = -1096
$ws.Range("N107").Value = -13838.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2937.327
$ws.Range("I31").Value = 1495.0714
$ws.Range("J31").Value = 8994.799999999999
$ws.Range("K31").Value = 1495.0714
$ws.Range("L31").Value = 8994.799999999999
$ws.Range("M31").Value = -1200.0714
$ws.Range("N31").Value = -9584.799999999999
$ws.Range("H34").Value = 2937.327
$ws.Range("I34").Value = 1495.0714
$ws.Range("J34").Value = 8994.799999999999
$ws.Range("K34").Value = 1495.0714
$ws.Range("L34").Value = 8994.799999999999
$ws.Range("M34").Value = -1293.0714
$ws.Range("N34").Value = -9398.799999999999
$ws.Range("H58").Value = 2581.5386
$ws.Range("I58").Value = 2194
$ws.Range("J58").Value = 3201.6
$ws.Range("K58").Value = 2194
$ws.Range("L58").Value = 3201.6
$ws.Range("M58").Value = -1991
$ws.Range("N58").Value = -3607.6
$ws.Range("H86").Value = 8760
$ws.Range("I86").Value = 6500
$ws.Range("J86").Value = 9325
$ws.Range("K86").Value = 6500
$ws.Range("L86").Value = 9325
$ws.Range("M86").Value = -5377
$ws.Range("N86").Value = -11571
$ws.Range("H89").Value = 8760
$ws.Range("I89").Value = 6500
$ws.Range("J89").Value = 9325
$ws.Range("K89").Value = 32500
$ws.Range("L89").Value = 46625
$ws.Range("M89").Value = -26884
$ws.Range("N89").Value = -57857
$ws.Range("H132").Value = 8339232.5
$ws.Range("I132").Value = 8775455
$ws.Range("J132").Value = 51002
$ws.Range("K132").Value = 26326365
$ws.Range("L132").Value = 153006
$ws.Range("M132").Value = -26323835
$ws.Range("N132").Value = -158066
$ws.Range("H134").Value = 1572.3793
$ws.Range("I134").Value = 1341.6041
$ws.Range("J134").Value = 2680.1
$ws.Range("K134").Value = 4024.8123
$ws.Range("L134").Value = 8040.299999999999
$ws.Range("M134").Value = -1489.8123
$ws.Range("N134").Value = -13110.3
$ws.Range("H136").Value = 2581.5386
$ws.Range("I136").Value = 2194
$ws.Range("J136").Value = 3201.6
$ws.Range("K136").Value = 6582
$ws.Range("L136").Value = 9604.799999999999
$ws.Range("M136").Value = -4032
$ws.Range("N136").Value = -14704.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1701.75
$ws.Range("I3").Value = 1701.75
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5105.25
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4993.25
$ws.Range("N3").ClearContents()
$ws.Range("H18").Value = 111794.11
$ws.Range("I18").Value = 143491.28
$ws.Range("J18").Value = 854
$ws.Range("K18").Value = 430473.84
$ws.Range("L18").Value = 2562
$ws.Range("M18").Value = -430304.84
$ws.Range("N18").Value = -2900
$ws.Range("H120").Value = 12015.25
$ws.Range("I120").Value = 8217.799999999999
$ws.Range("J120").Value = 18344.334
$ws.Range("K120").Value = 24653.4
$ws.Range("L120").Value = 55033.00199999999
$ws.Range("M120").Value = -19815.4
$ws.Range("N120").Value = -64709.00199999999
$ws.Range("H129").Value = 924.7692
$ws.Range("I129").Value = 924.7692
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2774.3076
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2225.6924
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 14099852
$ws.Range("I131").Value = 13973032
$ws.Range("J131").Value = 14150581
$ws.Range("K131").Value = 41919096
$ws.Range("L131").Value = 42451743
$ws.Range("M131").Value = -41914056
$ws.Range("N131").Value = -42461823
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 14996
$ws.Range("I44").Value = 14996
$ws.Range("K44").Value = 14996
$ws.Range("M44").Value = -14400
$ws.Range("H70").Value = 7238.095
$ws.Range("I70").Value = 6988
$ws.Range("J70").Value = 7571.5557
$ws.Range("K70").Value = 6988
$ws.Range("L70").Value = 7571.5557
$ws.Range("M70").Value = -6718
$ws.Range("N70").Value = -8111.5557
$ws.Range("H73").Value = 7238.095
$ws.Range("I73").Value = 6988
$ws.Range("J73").Value = 7571.5557
$ws.Range("K73").Value = 6988
$ws.Range("L73").Value = 7571.5557
$ws.Range("M73").Value = -6052
$ws.Range("N73").Value = -9443.555700000001
$ws.Range("H107").Value = 794.53845
$ws.Range("I107").Value = 599
$ws.Range("J107").Value = 853.2
$ws.Range("K107").Value = 599
$ws.Range("L107").Value = 853.2
$ws.Range("M107").Value = 1321
$ws.Range("N107").Value = -4693.2
$ws.Range("H126").Value = 4841.4116
$ws.Range("I126").Value = 2350
$ws.Range("J126").Value = 8400.571
$ws.Range("K126").Value = 7050
$ws.Range("L126").Value = 25201.713
$ws.Range("M126").Value = -4580
$ws.Range("N126").Value = -30141.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3315.7454
$ws.Range("I136").Value = 2365.658
$ws.Range("J136").Value = 5439.4707
$ws.Range("K136").Value = 7096.974
$ws.Range("L136").Value = 16318.4121
$ws.Range("M136").Value = -4546.974
$ws.Range("N136").Value = -21418.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3746404
$ws.Range("I132").Value = 4695611.5
$ws.Range("J132").Value = 2307.3333
$ws.Range("K132").Value = 14086834.5
$ws.Range("L132").Value = 6921.999899999999
$ws.Range("M132").Value = -14084304.5
$ws.Range("N132").Value = -11981.9999
$ws.Range("H136").Value = 7174.897
$ws.Range("I136").Value = 1296.9459
$ws.Range("J136").Value = 10799.634
$ws.Range("K136").Value = 3890.8377
$ws.Range("L136").Value = 32398.902
$ws.Range("M136").Value = -1340.8377
$ws.Range("N136").Value = -37498.902
